{"js": "// Update the multiplication problems in the first table of the document.\n// The table holds 5 \"blocks\" of problem rows (row indices 0, 4, 9, 14, 19),\n// each with 5 cells (columns 0-4), for a total of 25 problems that need\n// their text replaced in document order. Some original values repeat\n// (e.g. \"745\u00d76=\", \"840\u00d73=\") but map to different new values depending on\n// position, so cells are addressed positionally by (row, col) rather than\n// via a global text search/replace.\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"793\u00d76=\", newText: \"652\u00d76=\" },\n  { row: 0, col: 1, oldText: \"745\u00d76=\", newText: \"554\u00d77=\" },\n  { row: 0, col: 2, oldText: \"265\u00d78=\", newText: \"974\u00d77=\" },\n  { row: 0, col: 3, oldText: \"745\u00d76=\", newText: \"207\u00d76=\" },\n  { row: 0, col: 4, oldText: \"261\u00d79=\", newText: \"220\u00d72=\" },\n\n  { row: 4, col: 0, oldText: \"790\u00d74=\", newText: \"470\u00d79=\" },\n  { row: 4, col: 1, oldText: \"804\u00d78=\", newText: \"906\u00d72=\" },\n  { row: 4, col: 2, oldText: \"519\u00d74=\", newText: \"361\u00d74=\" },\n  { row: 4, col: 3, oldText: \"380\u00d73=\", newText: \"794\u00d75=\" },\n  { row: 4, col: 4, oldText: \"447\u00d75=\", newText: \"533\u00d78=\" },\n\n  { row: 9, col: 0, oldText: \"321\u00d77=\", newText: \"604\u00d77=\" },\n  { row: 9, col: 1, oldText: \"772\u00d77=\", newText: \"158\u00d79=\" },\n  { row: 9, col: 2, oldText: \"163\u00d74=\", newText: \"610\u00d76=\" },\n  { row: 9, col: 3, oldText: \"436\u00d78=\", newText: \"928\u00d74=\" },\n  { row: 9, col: 4, oldText: \"914\u00d77=\", newText: \"169\u00d77=\" },\n\n  { row: 14, col: 0, oldText: \"704\u00d74=\", newText: \"396\u00d73=\" },\n  { row: 14, col: 1, oldText: \"351\u00d74=\", newText: \"407\u00d74=\" },\n  { row: 14, col: 2, oldText: \"840\u00d73=\", newText: \"401\u00d78=\" },\n  { row: 14, col: 3, oldText: \"840\u00d73=\", newText: \"335\u00d76=\" },\n  { row: 14, col: 4, oldText: \"408\u00d74=\", newText: \"315\u00d78=\" },\n\n  { row: 19, col: 0, oldText: \"920\u00d77=\", newText: \"744\u00d73=\" },\n  { row: 19, col: 1, oldText: \"392\u00d75=\", newText: \"253\u00d74=\" },\n  { row: 19, col: 2, oldText: \"194\u00d73=\", newText: \"388\u00d73=\" },\n  { row: 19, col: 3, oldText: \"323\u00d75=\", newText: \"158\u00d79=\" },\n  { row: 19, col: 4, oldText: \"470\u00d78=\", newText: \"398\u00d77=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load current values first (sanity check), then write all the new values.\nconst cells = replacements.map(({ row, col }) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nreplacements.forEach(({ oldText, newText }, i) => {\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    console.log(\n      `Warning: expected \"${oldText}\" at cell #${i} but found \"${cell.value}\"`\n    );\n  }\n  cell.value = newText;\n});\nawait context.sync();\n", "ps1": "# Update the multiplication problems in the first table of the document.\n# The table holds 5 \"blocks\" of problem rows (1-based COM row indices\n# 1, 5, 10, 15, 20), each with 5 cells (columns 1-5), for a total of 25\n# problems whose text needs replacing, in document order. Some original\n# values repeat (e.g. \"745\u00d76=\", \"840\u00d73=\") but map to different new values\n# depending on position, so cells are addressed positionally via\n# Table.Cell(row, col) rather than via a global Find/Replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1;  Col = 1; Old = \"793\u00d76=\"; New = \"652\u00d76=\" },\n  @{ Row = 1;  Col = 2; Old = \"745\u00d76=\"; New = \"554\u00d77=\" },\n  @{ Row = 1;  Col = 3; Old = \"265\u00d78=\"; New = \"974\u00d77=\" },\n  @{ Row = 1;  Col = 4; Old = \"745\u00d76=\"; New = \"207\u00d76=\" },\n  @{ Row = 1;  Col = 5; Old = \"261\u00d79=\"; New = \"220\u00d72=\" },\n\n  @{ Row = 5;  Col = 1; Old = \"790\u00d74=\"; New = \"470\u00d79=\" },\n  @{ Row = 5;  Col = 2; Old = \"804\u00d78=\"; New = \"906\u00d72=\" },\n  @{ Row = 5;  Col = 3; Old = \"519\u00d74=\"; New = \"361\u00d74=\" },\n  @{ Row = 5;  Col = 4; Old = \"380\u00d73=\"; New = \"794\u00d75=\" },\n  @{ Row = 5;  Col = 5; Old = \"447\u00d75=\"; New = \"533\u00d78=\" },\n\n  @{ Row = 10; Col = 1; Old = \"321\u00d77=\"; New = \"604\u00d77=\" },\n  @{ Row = 10; Col = 2; Old = \"772\u00d77=\"; New = \"158\u00d79=\" },\n  @{ Row = 10; Col = 3; Old = \"163\u00d74=\"; New = \"610\u00d76=\" },\n  @{ Row = 10; Col = 4; Old = \"436\u00d78=\"; New = \"928\u00d74=\" },\n  @{ Row = 10; Col = 5; Old = \"914\u00d77=\"; New = \"169\u00d77=\" },\n\n  @{ Row = 15; Col = 1; Old = \"704\u00d74=\"; New = \"396\u00d73=\" },\n  @{ Row = 15; Col = 2; Old = \"351\u00d74=\"; New = \"407\u00d74=\" },\n  @{ Row = 15; Col = 3; Old = \"840\u00d73=\"; New = \"401\u00d78=\" },\n  @{ Row = 15; Col = 4; Old = \"840\u00d73=\"; New = \"335\u00d76=\" },\n  @{ Row = 15; Col = 5; Old = \"408\u00d74=\"; New = \"315\u00d78=\" },\n\n  @{ Row = 20; Col = 1; Old = \"920\u00d77=\"; New = \"744\u00d73=\" },\n  @{ Row = 20; Col = 2; Old = \"392\u00d75=\"; New = \"253\u00d74=\" },\n  @{ Row = 20; Col = 3; Old = \"194\u00d73=\"; New = \"388\u00d73=\" },\n  @{ Row = 20; Col = 4; Old = \"323\u00d75=\"; New = \"158\u00d79=\" },\n  @{ Row = 20; Col = 5; Old = \"470\u00d78=\"; New = \"398\u00d77=\" }\n)\n\nforeach ($r in $replacements) {\n  $cell = $t.Cell($r.Row, $r.Col)\n  $cellRange = $cell.Range\n  $current = $cellRange.Text.TrimEnd([char]7, \"`r\")\n  if ($current -ne $r.Old) {\n    Write-Output \"Warning: expected '$($r.Old)' at row $($r.Row) col $($r.Col) but found '$current'\"\n  }\n  $cellRange.Text = $r.New\n}\n"}
